# Generate Report for Handback
# Marks the 3c2fe154... and b6656fb4... documents as handed back (in sync
# with en-US) for both the zh-cn and de-de target languages, records the
# generated handback .xlf file names + handback timestamps, links the new
# "Latest Target File" cells back to the source doc, and widens a few
# columns on the per-language sheets (and the Overview roll-up) so the new
# content fits.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ColumnWidth (characters) -> stored <col width> is ColumnWidth + 5/6,
# rounded to the nearest 1/6. Solve backwards for the stored widths we need.
$wZhDe = 29.166666666666668   # -> stored ~29.98 (target column widths)
$wWide = 39.166666666666664   # -> stored 40

function Set-TargetLink {
    param(
        $ws,
        [string]$cell,
        [string]$mdName,
        [string]$url,
        [string]$xlfName,
        [string]$handbackDate
    )

    $ws.Range($cell).Value = $mdName
    $ws.Range($cell).Style = "Hyperlink"
    $ws.Range($cell).Font.Underline = 2
    $ws.Range($cell).Font.Color = 15570276
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status -> handed back, for both tracked docs
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (I) + Latest Handback File (J) + Latest Handback
# DateTime (K) for each row
$wsZh.Range("I2").Value = "3c2fe154-9406-4d69-919e-b0a487de73a4.md"
$wsZh.Range("J2").Value = "3c2fe154-9406-4d69-919e-b0a487de73a4.8587cb7082e3ed9c889fcdc6dcc21af3a91d85a5.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 14:51:58"

$wsZh.Range("I3").Value = "b6656fb4-fb24-47b5-83bf-9e01423676ee.md"
$wsZh.Range("J3").Value = "b6656fb4-fb24-47b5-83bf-9e01423676ee.afe84c42952cb26a99fc1a5cb005ba5491e883d7.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 14:51:58"

# Style the new "Latest Target File" links like the existing hyperlink cells
$wsZh.Range("I2").Style = "Hyperlink"
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Style = "Hyperlink"
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

# Rebuild hyperlinks in source-doc-order so rIds line up the way Excel
# would renumber them: A2, I2, A3, I3
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/3c2fe154-9406-4d69-919e-b0a487de73a4.md", "", "", "3c2fe154-9406-4d69-919e-b0a487de73a4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/3c2fe154-9406-4d69-919e-b0a487de73a4.md", "", "", "3c2fe154-9406-4d69-919e-b0a487de73a4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/b6656fb4-fb24-47b5-83bf-9e01423676ee.md", "", "", "b6656fb4-fb24-47b5-83bf-9e01423676ee.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/b6656fb4-fb24-47b5-83bf-9e01423676ee.md", "", "", "b6656fb4-fb24-47b5-83bf-9e01423676ee.md")

# Widen columns C (Status), I (Latest Target File), J (Latest Handback File)
$wsZh.Columns.Item(3).ColumnWidth = $wZhDe
$wsZh.Columns.Item(9).ColumnWidth = $wWide
$wsZh.Columns.Item(10).ColumnWidth = $wWide

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "3c2fe154-9406-4d69-919e-b0a487de73a4.md"
$wsDe.Range("J2").Value = "3c2fe154-9406-4d69-919e-b0a487de73a4.8587cb7082e3ed9c889fcdc6dcc21af3a91d85a5.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 14:52:11"

$wsDe.Range("I3").Value = "b6656fb4-fb24-47b5-83bf-9e01423676ee.md"
$wsDe.Range("J3").Value = "b6656fb4-fb24-47b5-83bf-9e01423676ee.afe84c42952cb26a99fc1a5cb005ba5491e883d7.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 14:52:11"

$wsDe.Range("I2").Style = "Hyperlink"
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Style = "Hyperlink"
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/3c2fe154-9406-4d69-919e-b0a487de73a4.md", "", "", "3c2fe154-9406-4d69-919e-b0a487de73a4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/3c2fe154-9406-4d69-919e-b0a487de73a4.md", "", "", "3c2fe154-9406-4d69-919e-b0a487de73a4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/b6656fb4-fb24-47b5-83bf-9e01423676ee.md", "", "", "b6656fb4-fb24-47b5-83bf-9e01423676ee.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7593befbc9b57e17ec98bfb2ec0f8d903e407f2/e2e/b6656fb4-fb24-47b5-83bf-9e01423676ee.md", "", "", "b6656fb4-fb24-47b5-83bf-9e01423676ee.md")

$wsDe.Columns.Item(3).ColumnWidth = $wZhDe
$wsDe.Columns.Item(9).ColumnWidth = $wWide
$wsDe.Columns.Item(10).ColumnWidth = $wWide

# ---------------------------------------------------------------------
# Overview roll-up sheet: widen the zh-cn / de-de columns (E, F) to match
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $wZhDe
$wsOverview.Columns.Item(6).ColumnWidth = $wZhDe
